# Age_OM.xlsx restructuring:
#  - relabel the "w" header as "w (cm/year)" (column E) and duplicate it onto
#    the new column F header
#  - add a second "at Xcm" age-calculation block in column J ("at 5cm")
#  - add a new "t=10years" column M that derives depth from the sedimentation
#    rate and a fixed 10-year duration
#  - update the saved selection to I21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 headers ---------------------------------------------------------
# "w" -> "w (cm/year)"
$ws.Range("E3").Value = "w (cm/year)"
# New column F reuses the same header text as E3
$ws.Range("F3").Value = "w (cm/year)"
# New header columns for the additional age computations
$ws.Range("J3").Value = "at 5cm"
$ws.Range("M3").Value = "t=10years"

# --- New formulas for rows 4-5 ---------------------------------------------
# Column J: age (days) at 5cm depth, same pattern as existing column I (50cm)
$ws.Range("J4").Formula = "=5/F4"
$ws.Range("J5").Formula = "=5/F5"

# Column M: depth (cm) reached after a fixed 10 years
$ws.Range("M4").Formula = "=10*F4"
$ws.Range("M5").Formula = "=10*F5"

# The new J4:J5 cells pick up a distinct (re-applied "Normal") cell style,
# separate from the default style used everywhere else
$ws.Range("J4:J5").Style = "Normal"

# --- Restore the saved selection -------------------------------------------
$ws.Range("I21").Select()
